$d = $word.ActiveDocument

# 1) Programa section (Portuguese), first block with numbered list 1-4
$old1 = "1) Apresentação da disciplina e grade curricular do curso de Engenharia Química na EEL.2) A Engenharia Química e a Sociedade - Histórico e Objetivos – A Eng. Química na Indústria, no Ensino e Pesquisa - Atribuições e Campo de Atuação do Engenheiro Químico - O Mercado de Trabalho - Órgãos que Regulamentam a profissão do Engenheiro Químico. Associações de classe.3) Definições básicas e Operações Unitárias na Eng. Química – Interpretação de fluxogramas de processo - As Indústrias Químicas e os Processos Industriais – Principais dimensões e unidades na engenharia química e suas conversões.4) Palestras e Visitas Técnicas a EEL e/ou a indústrias químicas para conhecimento de processos."
$new1 = "1) Apresentação da disciplina e grade curricular do curso de Engenharia Química na EEL.^l2) A Engenharia Química e a Sociedade - Histórico e Objetivos – A Eng. Química na Indústria, no Ensino e Pesquisa - Atribuições e Campo de Atuação do Engenheiro Químico - O Mercado de Trabalho - Órgãos que Regulamentam a profissão do Engenheiro Químico. Associações de classe.^l3) Definições básicas e Operações Unitárias na Eng. Química – Interpretação de fluxogramas de processo - As Indústrias Químicas e os Processos Industriais – Principais dimensões e unidades na engenharia química e suas conversões.^l4) Palestras e Visitas Técnicas a EEL e/ou a indústrias químicas para conhecimento de processos."
$f1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Host "Replace 1:" $f1

# 2) Média Final line
$old2 = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) / 3Média final mínima de aprovação = 5,0"
$new2 = "Média Final = (Prova1 + Prova2 + Nota de Trabalho) / 3^lMédia final mínima de aprovação = 5,0"
$f2 = $d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)
Write-Host "Replace 2:" $f2

# 3) Bibliografia list (Portuguese)
$old3 = "1) INDÚSTRIA QUÍMICA – RISCOS E OPORTUNIDADES, Pedro Wongtschowski, Ed. Edgar Blucher, 2002.2) INDÚSTRIAS DE PROCESSOS QUÍMICOS SHREVE, R. Norris Jr. ; Joseph A. Brink Ed. LTC, 1980.3) PRINCÍPIOS ELEMENTARES DE PROCESSOS QUÍMICOS, Felder, R.M.; Roussau, R.W. , 2005.4) ENGENHARIA QUÍMICA - PRINCÍPIOS E CÁLCULOS, Himmelblau, D.M. - Riggs, J.B.  Ed. LTC,  20065) INTRODUÇÃO A ENGENHARIA QUÍMICA BRASIL, Nilo Índio Ed. Interciência , 2013.6) ENGENHARIA QUÍMICA, Cremasco, Marco Aurélio, Ed. Edgard Blucher, 2005."
$new3 = "1) INDÚSTRIA QUÍMICA – RISCOS E OPORTUNIDADES, Pedro Wongtschowski, Ed. Edgar Blucher, 2002.^l2) INDÚSTRIAS DE PROCESSOS QUÍMICOS SHREVE, R. Norris Jr. ; Joseph A. Brink Ed. LTC, 1980.^l3) PRINCÍPIOS ELEMENTARES DE PROCESSOS QUÍMICOS, Felder, R.M.; Roussau, R.W. , 2005.^l4) ENGENHARIA QUÍMICA - PRINCÍPIOS E CÁLCULOS, Himmelblau, D.M. - Riggs, J.B.  Ed. LTC,  2006^l5) INTRODUÇÃO A ENGENHARIA QUÍMICA BRASIL, Nilo Índio Ed. Interciência , 2013.^l6) ENGENHARIA QUÍMICA, Cremasco, Marco Aurélio, Ed. Edgard Blucher, 2005."
$f3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Host "Replace 3:" $f3

# 4) Programa list (English), italic run
$old4 = "1) Presentation of the discipline and curriculum guides of the Chemical Engineering course at EEL.2) Chemical Engineering and Society - History and Objectives - Chemical Engineering in Industry, Education and Research - Attributions and areas of activity of the Chemical Engineer - The Labor Market - Regulating the profession of the Chemical Engineer. Class associations.3) Basic Definitions and Unit Operations in Chemical Engineering - Interpretation of Process Flowcharts - Chemical Industries and Industrial Processes - Major dimensions and units in chemical engineering and their conversions.4) Lectures and Technical Visits to EEL and / or the chemical industries for process knowledge."
$new4 = "1) Presentation of the discipline and curriculum guides of the Chemical Engineering course at EEL.^l2) Chemical Engineering and Society - History and Objectives - Chemical Engineering in Industry, Education and Research - Attributions and areas of activity of the Chemical Engineer - The Labor Market - Regulating the profession of the Chemical Engineer. Class associations.^l3) Basic Definitions and Unit Operations in Chemical Engineering - Interpretation of Process Flowcharts - Chemical Industries and Industrial Processes - Major dimensions and units in chemical engineering and their conversions.^l4) Lectures and Technical Visits to EEL and / or the chemical industries for process knowledge."
$f4 = $d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Host "Replace 4:" $f4

# 5) Critério paragraph - long text with bullet points
$old5 = "Os objetivos esperados para colaboração na formação dos discentes englobam fortalecer a percepção da carreira em formação, bem como motivá-los aos desafios a que estarão direcionados no decorrer no curso.3Quanto aos participantes da Indústria espera-se estreitar as relações Escola/Empresa além de permitir a apresentação de opiniões concretas sobre o que a Indústria espera da formação de nossos alunos, além disso, identificar o que as indústrias podem colaborar na formação acadêmica dos alunos. Essa colaboração engloba diversos aspectos que beneficiam ambas as partes, como pode ser observado a seguir.Para os discentes:• Fortalecimento da percepção na carreira: interação com profissionais da indústria proporciona aos discentes uma visão mais clara e realista sobre as expectativas e desafios da carreira que estão buscando. Isso ajuda a fortalecer sua motivação e engajamento com o curso de Engenharia Química.• Preparação para desafios futuros: Ao participar de atividades práticas e colaborativas com profissionais da indústria, os discentes são expostos aos desafios reais enfrentados no campo, preparando-os melhor para enfrentar esses desafios durante e após a conclusão do curso.Para os participantes da indústria:• Estreitamento das relações escola/ empresa: A colaboração com a universidade proporciona uma oportunidade valiosa para as indústrias se envolverem ativamente na formação dos futuros profissionais da área. Isso fortalece os laços entre a academia e o setor industrial, criando uma relação de parceria e colaboração mútua.• Apresentação de expectativas e necessidades das indústrias: Os participantes da indústria têm a oportunidade de expressar suas opiniões e expectativas em relação à formação dos alunos, destacando as habilidades e conhecimentos que são valorizados no mercado de trabalho. Isso ajuda a alinhar o currículo acadêmico com as demandas da indústria.• Identificação de oportunidades de colaboração: A colaboração entre a universidade e a indústria pode abrir portas para futuras parcerias e projetos conjuntos de pesquisa e desenvolvimento. Isso permite que as indústrias contribuam ativamente para a formação acadêmica dos alunos, oferecendo oportunidades de estágio, projetos de conclusão de curso e outras experiências práticas relevantes. A colaboração entre os discentes e os participantes da indústria beneficia ambas as partes, fortalecendo a formação dos alunos e promovendo uma maior integração entre a academia e o setor industrial."
$new5 = "Os objetivos esperados para colaboração na formação dos discentes englobam fortalecer a percepção da carreira em formação, bem como motivá-los aos desafios a que estarão direcionados no decorrer no curso.^l3^lQuanto aos participantes da Indústria espera-se estreitar as relações Escola/Empresa além de permitir a apresentação de opiniões concretas sobre o que a Indústria espera da formação de nossos alunos, além disso, identificar o que as indústrias podem colaborar na formação acadêmica dos alunos. Essa colaboração engloba diversos aspectos que beneficiam ambas as partes, como pode ser observado a seguir.^lPara os discentes:^l• Fortalecimento da percepção na carreira: interação com profissionais da indústria proporciona aos discentes uma visão mais clara e realista sobre as expectativas e desafios da carreira que estão buscando. Isso ajuda a fortalecer sua motivação e engajamento com o curso de Engenharia Química.^l• Preparação para desafios futuros: Ao participar de atividades práticas e colaborativas com profissionais da indústria, os discentes são expostos aos desafios reais enfrentados no campo, preparando-os melhor para enfrentar esses desafios durante e após a conclusão do curso.^lPara os participantes da indústria:^l• Estreitamento das relações escola/ empresa: A colaboração com a universidade proporciona uma oportunidade valiosa para as indústrias se envolverem ativamente na formação dos futuros profissionais da área. Isso fortalece os laços entre a academia e o setor industrial, criando uma relação de parceria e colaboração mútua.^l• Apresentação de expectativas e necessidades das indústrias: Os participantes da indústria têm a oportunidade de expressar suas opiniões e expectativas em relação à formação dos alunos, destacando as habilidades e conhecimentos que são valorizados no mercado de trabalho. Isso ajuda a alinhar o currículo acadêmico com as demandas da indústria.^l• Identificação de oportunidades de colaboração: A colaboração entre a universidade e a indústria pode abrir portas para futuras parcerias e projetos conjuntos de pesquisa e desenvolvimento. Isso permite que as indústrias contribuam ativamente para a formação acadêmica dos alunos, oferecendo oportunidades de estágio, projetos de conclusão de curso e outras experiências práticas relevantes. A colaboração entre os discentes e os participantes da indústria beneficia ambas as partes, fortalecendo a formação dos alunos e promovendo uma maior integração entre a academia e o setor industrial."
$f5 = $d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)
Write-Host "Replace 5:" $f5
